$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Record", "Balanço Geral", "Trânsito", "2025-04-10T18:10", "Neutro", "Carro cai dentro de canal às margens da Estrada do Carvão: ninguém se feriu. Repórter *ao vivo*. Ainda não se sabe as causas do acidente. Repórter falou sobre a estrada ser uma via alternativa para caminhões acima de 4 eixos, que estão proibidos de circular dentro do município. Apresentador lembrou de outro caso recente em que carro de um carroceiro caiu no canal e ele acabou morrendo.  "),
    @("Record", "RJ Record", "Trânsito", "2025-04-10T18:11", "Neutro", "Campos: ciclista fica ferido após ser atropelado por carro na RJ-216. Imagens da vítima recebendo os primeiros atendimentos e aguardando socorro do Corpo de Bombeiros. *nota coberta*"),
    @("Record", "RJ Record", "Trânsito", "2025-04-10T18:19", "Neutro", "Notícias em destaque na região: colisão entre duas motos no cruzamento entre as ruas Cora de Alvarenga e Epitácio Pessoa. Colisão aconteceu ontem à tarde. Acidentes são frequentes no local."),
    @("Record", "RJ Record", "Infraestrutura", "2025-04-10T18:23", "Negativo", "Mudança climática. Chuva forte causa alagamentos em cidades das regiões Norte e Lagos. Repórter *ao vivo*. Em Campos, ruas ficaram alagadas. Chuva rápida e intensa. Chuva confirmou previsão do INMET. "),
    @("Record", "RJ Record", "Educação", "2025-04-10T18:26", "Positivo", "Continua o recadastramento ara usuários do transporte universitário. Repórter *ao vivo*. Prazo é até 30 de abril. Vídeo com a subsecretária, Cátia Melo. Tema foi abordado no Balanço Geral. ")
)

$startRow = 119
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
